$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header field updates ---
# C3 holds a numeric-looking identifier that must remain text, matching the
# original cell's string type (not auto-converted to a number).
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "2556845"
$ws.Range("C3").ClearFormats()

$ws.Range("C4").Value = "DISEÑO E INTEGRACION DE MULTIMEDIA"
$ws.Range("C6").Value = 44760
$ws.Range("C7").Value = 45124

# --- Existing rows 11-15: update data ---
$data = @(
    @("CC", 1005178211, "JUAN CAMILO", "DELGADO CARRASCAL", "CERTIFICADO"),
    @("CC", 1005181992, "JOHAN", "VARGAS CALDERIN", "POR CERTIFICAR"),
    @("CC", 1005184329, "WILLIAM ANDRES", "LOPEZ RIOS", "CANCELADO"),
    @("CC", 1005185919, "SEBASTIAN", "PERTUZ SAMPAYO", "CERTIFICADO"),
    @("CC", 1005220651, "BRAYAN EDUARDO", "BADILLO HERRERA", "CERTIFICADO"),
    @("CC", 1005239745, "SARAY DUVIANA", "UNRIZA JAIMES", "CERTIFICADO"),
    @("CC", 1005241421, "CLARA LUCIA", "RUIZ MONSALVE", "RETIRO VOLUNTARIO"),
    @("CC", 1043962939, "DANNA KAROLAY", "RESTREPO SOSA", "CERTIFICADO"),
    @("CC", 1048457729, "DAYANA", "URRUCHURTU NIÑO", "CERTIFICADO"),
    @("TI", 1049019898, "KAREN YURLEIDY", "MARIN VARGAS", "RETIRO VOLUNTARIO"),
    @("CC", 1087985197, "GISELL MARIANA", "MARIN LARROTA", "CERTIFICADO"),
    @("CC", 1096184002, "DANIELA", "ROJAS BOTELLO", "CERTIFICADO"),
    @("CC", 1096186262, "KEVIN ANDRES", "PARADA SUAREZ", "RETIRO VOLUNTARIO"),
    @("CC", 1096189477, "KAMILA", "QUINTERO CARREÑO", "CERTIFICADO"),
    @("CC", 1097183074, "MARIA JOSE", "ORTIZ GUIZA", "CERTIFICADO"),
    @("CC", 1144182405, "CAROLAIN", "ABANIS PEREZ", "CERTIFICADO"),
    @("CC", 63469380, "VIDA EMPERATRIZ", "SANTOS YAIN", "CERTIFICADO")
)

$startRow = 11
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
}
